# Interview-Guide.xlsx — add "Bubble sort" and "Stack" rows to the guide.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Bubble sort (Sorting algorithm sub-topic)
$ws.Range("F7").Value = "Sorting algorithm "
$ws.Range("G7").Value = "Bubble sort"
$ws.Range("H7").Value = "pair are compared and sawped for sorting"
$ws.Range("I7").Value = "O(n^2)"
$ws.Range("J7").Value = "space O(1)"

# Row 8: Stack (Data Structure sub-topic)
$ws.Range("F8").Value = "Data Structure "
$ws.Range("G8").Value = "Stack"
$ws.Range("H8").Value = "First in Last out`nLast in first out"
$ws.Range("I8").Value = "O(n) for search and access`nO(1) for insertion and deletion"
$ws.Range("J8").Value = "O(n) worst"
$ws.Range("K8").Value = "For balancing eqaution it can be used."

# Row 8 needs the taller row height used by the other multi-line rows.
$ws.Rows.Item(8).RowHeight = 60

# Restore the scrolled-down selection state (user had scrolled to row 8).
$ws.Range("K8").Select()
